$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs / Vwf / Tnfrsf11b / FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vwf"
$ws.Range("C2").Value = "Tnfrsf11b"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 33.00563866666667
$ws.Range("H2").Value = 99.01691600000001
$ws.Range("I2").Value = 0.9169150302490913
$ws.Range("J2").Value = 0.9169150302490912
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.776574666666666
$ws.Range("N2").Value = 11.329724
$ws.Range("O2").Value = 0.9855052394405499
$ws.Range("P2").Value = 0.9855052394405499
$ws.Range("Q2").Value = 124.6482588456871
$ws.Range("R2").Value = 1121.834329611184
$ws.Range("S2").Value = 0.9036245664322697
$ws.Range("T2").Value = 0.9036245664322696

# Row 3: ECs / Vwf / Tnfrsf11b / sCs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vwf"
$ws.Range("C3").Value = "Tnfrsf11b"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 33.00563866666667
$ws.Range("H3").Value = 99.01691600000001
$ws.Range("I3").Value = 0.9169150302490913
$ws.Range("J3").Value = 0.9169150302490912
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.05554566666666667
$ws.Range("N3").Value = 0.166637
$ws.Range("O3").Value = 0.01449476055945007
$ws.Range("P3").Value = 0.01449476055945008
$ws.Range("Q3").Value = 1.833320203499111
$ws.Range("R3").Value = 16.499881831492
$ws.Range("S3").Value = 0.0132904638168215
$ws.Range("T3").Value = 0.0132904638168215

# Row 4: FAPs / Vwf / Tnfrsf11b / FAPs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Vwf"
$ws.Range("C4").Value = "Tnfrsf11b"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.07465466666666666
$ws.Range("H4").Value = 0.223964
$ws.Range("I4").Value = 0.002073948231580021
$ws.Range("J4").Value = 0.00207394823158002
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.776574666666666
$ws.Range("N4").Value = 11.329724
$ws.Range("O4").Value = 0.9855052394405499
$ws.Range("P4").Value = 0.9855052394405499
$ws.Range("Q4").Value = 0.2819389228817777
$ws.Range("R4").Value = 2.537450305936
$ws.Range("S4").Value = 0.002043886848550574
$ws.Range("T4").Value = 0.002043886848550573

# Row 5: FAPs / Vwf / Tnfrsf11b / sCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Vwf"
$ws.Range("C5").Value = "Tnfrsf11b"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.07465466666666666
$ws.Range("H5").Value = 0.223964
$ws.Range("I5").Value = 0.002073948231580021
$ws.Range("J5").Value = 0.00207394823158002
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.05554566666666667
$ws.Range("N5").Value = 0.166637
$ws.Range("O5").Value = 0.01449476055945007
$ws.Range("P5").Value = 0.01449476055945008
$ws.Range("Q5").Value = 0.004146743229777777
$ws.Range("R5").Value = 0.037320689068
$ws.Range("S5").Value = 0.00003006138302944731
$ws.Range("T5").Value = 0.00003006138302944731

# Row 6: M2 / Vwf / Tnfrsf11b / FAPs
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Vwf"
$ws.Range("C6").Value = "Tnfrsf11b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.576299333333333
$ws.Range("H6").Value = 7.728898
$ws.Range("I6").Value = 0.07157103078692272
$ws.Range("J6").Value = 0.0715710307869227
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.776574666666666
$ws.Range("N6").Value = 11.329724
$ws.Range("O6").Value = 0.9855052394405499
$ws.Range("P6").Value = 0.9855052394405499
$ws.Range("Q6").Value = 9.729586796016889
$ws.Range("R6").Value = 87.56628116415199
$ws.Range("S6").Value = 0.07053362583267324
$ws.Range("T6").Value = 0.07053362583267322

# Row 7: M2 / Vwf / Tnfrsf11b / sCs
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Vwf"
$ws.Range("C7").Value = "Tnfrsf11b"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.576299333333333
$ws.Range("H7").Value = 7.728898
$ws.Range("I7").Value = 0.07157103078692272
$ws.Range("J7").Value = 0.0715710307869227
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.05554566666666667
$ws.Range("N7").Value = 0.166637
$ws.Range("O7").Value = 0.01449476055945007
$ws.Range("P7").Value = 0.01449476055945008
$ws.Range("Q7").Value = 0.1431022640028889
$ws.Range("R7").Value = 1.287920376026
$ws.Range("S7").Value = 0.001037404954249474
$ws.Range("T7").Value = 0.001037404954249474

# Row 8: sCs / Vwf / Tnfrsf11b / FAPs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Vwf"
$ws.Range("C8").Value = "Tnfrsf11b"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.3398056666666667
$ws.Range("H8").Value = 1.019417
$ws.Range("I8").Value = 0.009439990732406145
$ws.Range("J8").Value = 0.009439990732406145
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.776574666666666
$ws.Range("N8").Value = 11.329724
$ws.Range("O8").Value = 0.9855052394405499
$ws.Range("P8").Value = 0.9855052394405499
$ws.Range("Q8").Value = 1.283301472323111
$ws.Range("R8").Value = 11.549713250908
$ws.Range("S8").Value = 0.009303160327056491
$ws.Range("T8").Value = 0.009303160327056491

# Row 9: sCs / Vwf / Tnfrsf11b / sCs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Vwf"
$ws.Range("C9").Value = "Tnfrsf11b"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.3398056666666667
$ws.Range("H9").Value = 1.019417
$ws.Range("I9").Value = 0.009439990732406145
$ws.Range("J9").Value = 0.009439990732406145
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.05554566666666667
$ws.Range("N9").Value = 0.166637
$ws.Range("O9").Value = 0.01449476055945007
$ws.Range("P9").Value = 0.01449476055945008
$ws.Range("Q9").Value = 0.01887473229211111
$ws.Range("R9").Value = 0.169872590629
$ws.Range("S9").Value = 0.0001368304053496548
$ws.Range("T9").Value = 0.0001368304053496548
